# This workbook edit ("Loan RBI, Variable Instalments"):
#  - inserts a new (blank) column N on the "Repayment Schedule" sheet,
#    pushing the previous N/O/P columns (Late / Heading / Amount) to O/P/Q
#  - widens columns K:L slightly and gives the newly inserted column N
#    the same width as column M (mirrors what Excel does automatically
#    when a column is inserted)
#  - makes "Repayment Schedule" the active/selected sheet (was "Summary")
#    with H17 as the selected cell

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$schedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14th column).
# Everything that used to live in N/O/P shifts right to O/P/Q.
$schedule.Columns("N").Insert()

# Match the look-and-feel Excel applies automatically on column insert:
# the new blank column inherits the width of the column to its left (M),
# and the K:L columns end up very slightly wider.
$schedule.Columns("N").ColumnWidth = $schedule.Columns("M").ColumnWidth
$schedule.Columns("K:L").ColumnWidth = 7.4

# Make "Repayment Schedule" the active sheet/tab (previously "Summary" was).
$schedule.Activate()
$schedule.Range("H17").Select()
